# Auto-generated Word COM-interop script
# Applies the adc_sbb_within_100 daily-update edit:
#   - bumps the heading date by one day
#   - replaces each of the 100 arithmetic problems in the table,
#     cell-by-cell (in document/row-major order) so that duplicate
#     old expressions (e.g. "24+19=", "81-6=") each get the correct,
#     position-specific replacement.
#
# Notes on Find.Execute args used below:
#   MatchWholeWord=$true  -- "=" is not a word boundary, so without this
#                            "9+44=" would match inside "49+44=" etc.
#   Wrap=0 (wdFindStop)   -- do not let the search spill past the cell.
#   Replace=1 (wdReplaceOne) -- wdReplaceAll(2) replaces every matching
#                            occurrence in the whole story, which breaks
#                            cells that happen to share identical text
#                            (e.g. the two "81-6=" / "24+19=" cells).

$d = $word.ActiveDocument

# --- Heading date ---
$d.Content.Find.Execute("2024-02-12 Monday", $true, $true, $false, $false, $false, $true, 0, $false, "2024-02-13 Tuesday", 1) | Out-Null

# --- Table of 100 arithmetic problems (20 rows x 5 columns) ---
$pairs = @(
    @{old="9+44="; new="36+15="},
    @{old="34+9="; new="68+24="},
    @{old="72-9="; new="71-35="},
    @{old="18+24="; new="6+75="},
    @{old="14+39="; new="90-83="},
    @{old="63-6="; new="89+7="},
    @{old="24+19="; new="57-19="},
    @{old="57+7="; new="45-28="},
    @{old="29+43="; new="50-26="},
    @{old="18+8="; new="79+15="},
    @{old="19+55="; new="40-3="},
    @{old="70-64="; new="91-12="},
    @{old="49+15="; new="49+36="},
    @{old="81-6="; new="39+17="},
    @{old="66-57="; new="84-55="},
    @{old="74+9="; new="80-7="},
    @{old="49+44="; new="40-13="},
    @{old="19+62="; new="13+58="},
    @{old="36+49="; new="4+67="},
    @{old="38+14="; new="80-69="},
    @{old="32-27="; new="37+15="},
    @{old="97-38="; new="76-29="},
    @{old="90-43="; new="71-14="},
    @{old="22-16="; new="38+6="},
    @{old="60-58="; new="52-34="},
    @{old="14+48="; new="7+87="},
    @{old="80-3="; new="6+65="},
    @{old="92-86="; new="22-15="},
    @{old="7+48="; new="28+8="},
    @{old="34+19="; new="94-75="},
    @{old="74-16="; new="66-18="},
    @{old="59+9="; new="94-7="},
    @{old="52-38="; new="64+28="},
    @{old="18+15="; new="5+58="},
    @{old="75-8="; new="46+27="},
    @{old="71-18="; new="37+39="},
    @{old="68+15="; new="39+17="},
    @{old="81-6="; new="56+27="},
    @{old="26+8="; new="91-66="},
    @{old="56+39="; new="39+15="},
    @{old="35+48="; new="76-38="},
    @{old="24+19="; new="44+17="},
    @{old="26+68="; new="3+58="},
    @{old="2+29="; new="4+77="},
    @{old="22+49="; new="62-8="},
    @{old="49+22="; new="29+12="},
    @{old="63-29="; new="59+14="},
    @{old="92-35="; new="43-24="},
    @{old="3+89="; new="47+4="},
    @{old="77-9="; new="51-33="},
    @{old="42-19="; new="67-59="},
    @{old="17+65="; new="64-7="},
    @{old="50-6="; new="9+72="},
    @{old="24+68="; new="72-63="},
    @{old="26+27="; new="4+17="},
    @{old="80-66="; new="60-33="},
    @{old="81-8="; new="17+78="},
    @{old="72-26="; new="17-8="},
    @{old="92-85="; new="69+29="},
    @{old="49+8="; new="67-39="},
    @{old="61-29="; new="41-37="},
    @{old="82-54="; new="91-7="},
    @{old="30-18="; new="36-7="},
    @{old="52-5="; new="26-18="},
    @{old="77-48="; new="80-44="},
    @{old="82-58="; new="17+35="},
    @{old="81-54="; new="43+38="},
    @{old="32-15="; new="63-24="},
    @{old="38+4="; new="41-17="},
    @{old="65-6="; new="7+6="},
    @{old="87+9="; new="12+9="},
    @{old="8+83="; new="70-44="},
    @{old="6+29="; new="40-18="},
    @{old="46-28="; new="16+37="},
    @{old="8+77="; new="88+8="},
    @{old="63-9="; new="70-24="},
    @{old="64-27="; new="6+65="},
    @{old="64+7="; new="70-31="},
    @{old="90-47="; new="56+37="},
    @{old="22-19="; new="72-66="},
    @{old="44+27="; new="84-17="},
    @{old="28+17="; new="91-29="},
    @{old="46+8="; new="91-73="},
    @{old="81-78="; new="51-29="},
    @{old="25+9="; new="19+15="},
    @{old="75-69="; new="14-9="},
    @{old="80-45="; new="40-14="},
    @{old="57-48="; new="53-44="},
    @{old="93-39="; new="28+67="},
    @{old="59+29="; new="27+4="},
    @{old="94-79="; new="22-18="},
    @{old="48+44="; new="7+78="},
    @{old="49+27="; new="92-38="},
    @{old="68+23="; new="72-4="},
    @{old="60-14="; new="73-15="},
    @{old="61-28="; new="8+13="},
    @{old="7+15="; new="79+6="},
    @{old="39+57="; new="86+6="},
    @{old="9+2="; new="17+8="},
    @{old="44-38="; new="27+68="}
)

$t = $d.Tables.Item(1)
$idx = 0
for ($row = 1; $row -le $t.Rows.Count; $row++) {
    for ($col = 1; $col -le $t.Columns.Count; $col++) {
        $pair = $pairs[$idx]
        $cell = $t.Cell($row, $col)
        $cell.Range.Find.Execute($pair.old, $true, $true, $false, $false, $false, $true, 0, $false, $pair.new, 1) | Out-Null
        $idx++
    }
}

Write-Host "Replaced $idx cells"